# "Coming soon" doc: strip the placeholder text, leaving a single empty
# paragraph behind (the paragraph mark itself is preserved).
$d = $word.ActiveDocument

# Select the entire body content and clear it. This removes every run
# (the "Coming soon:" / " How to " / "digitally " / "implement a
# biquadratic notch filter" runs) but keeps the trailing paragraph mark,
# so the document ends up with exactly one empty paragraph.
$d.Content.Delete()
